$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.398441
$ws.Range("H2").Value = 0.796882
$ws.Range("I2").Value = 0.08945363909080989
$ws.Range("J2").Value = 0.06146862341190577
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.4984935
$ws.Range("N2").Value = 12.996987
$ws.Range("O2").Value = 0.136717006960842
$ws.Range("P2").Value = 0.1171455833533661
$ws.Range("Q2").Value = 2.5892662486335
$ws.Range("R2").Value = 10.357064994534
$ws.Range("S2").Value = 0.01222983379825091
$ws.Range("T2").Value = 0.007200777747516081

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.398441
$ws.Range("H3").Value = 0.796882
$ws.Range("I3").Value = 0.08945363909080989
$ws.Range("J3").Value = 0.06146862341190577
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.026005
$ws.Range("N3").Value = 21.078015
$ws.Range("O3").Value = 0.1478149319518302
$ws.Range("P3").Value = 0.1899822138089391
$ws.Range("Q3").Value = 2.799448458205
$ws.Range("R3").Value = 16.79669074923
$ws.Range("S3").Value = 0.01322258357505164
$ws.Range("T3").Value = 0.01167794515558184

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.398441
$ws.Range("H4").Value = 0.796882
$ws.Range("I4").Value = 0.08945363909080989
$ws.Range("J4").Value = 0.06146862341190577
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.048527
$ws.Range("N4").Value = 6.145581
$ws.Range("O4").Value = 0.04309744714193724
$ws.Range("P4").Value = 0.05539188977340388
$ws.Range("Q4").Value = 0.8162171464069999
$ws.Range("R4").Value = 4.897302878442
$ws.Range("S4").Value = 0.003855223482370111
$ws.Range("T4").Value = 0.003404863212555157

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.398441
$ws.Range("H5").Value = 0.796882
$ws.Range("I5").Value = 0.08945363909080989
$ws.Range("J5").Value = 0.06146862341190577
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.046790666666667
$ws.Range("N5").Value = 9.140372000000001
$ws.Range("O5").Value = 0.0640991794148744
$ws.Range("P5").Value = 0.08238480272441405
$ws.Range("Q5").Value = 1.213966320017333
$ws.Range("R5").Value = 7.283797920104001
$ws.Range("S5").Value = 0.005733904861395246
$ws.Range("T5").Value = 0.005064080413531156

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.398441
$ws.Range("H6").Value = 0.796882
$ws.Range("I6").Value = 0.08945363909080989
$ws.Range("J6").Value = 0.06146862341190577
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.761096
$ws.Range("N6").Value = 11.283288
$ws.Range("O6").Value = 0.0791269219569728
$ws.Range("P6").Value = 0.1016995211970309
$ws.Range("Q6").Value = 1.498574851336
$ws.Range("R6").Value = 8.991449108015999
$ws.Range("S6").Value = 0.007078191119105725
$ws.Range("T6").Value = 0.006251329569631424

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.398441
$ws.Range("H7").Value = 0.796882
$ws.Range("I7").Value = 0.08945363909080989
$ws.Range("J7").Value = 0.06146862341190577
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 25.151532
$ws.Range("N7").Value = 50.303064
$ws.Range("O7").Value = 0.5291445125735434
$ws.Range("P7").Value = 0.453395989142846
$ws.Range("Q7").Value = 10.021401561612
$ws.Range("R7").Value = 40.08560624644799
$ws.Range("S7").Value = 0.04733390225463627
$ws.Range("T7").Value = 0.02786962731309012

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.055721
$ws.Range("H8").Value = 12.167163
$ws.Range("I8").Value = 0.9105463609091901
$ws.Range("J8").Value = 0.9385313765880943
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.4984935
$ws.Range("N8").Value = 12.996987
$ws.Range("O8").Value = 0.136717006960842
$ws.Range("P8").Value = 0.1171455833533661
$ws.Range("Q8").Value = 26.3560765563135
$ws.Range("R8").Value = 158.136459337881
$ws.Range("S8").Value = 0.1244871731625911
$ws.Range("T8").Value = 0.10994480560585

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.055721
$ws.Range("H9").Value = 12.167163
$ws.Range("I9").Value = 0.9105463609091901
$ws.Range("J9").Value = 0.9385313765880943
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.026005
$ws.Range("N9").Value = 21.078015
$ws.Range("O9").Value = 0.1478149319518302
$ws.Range("P9").Value = 0.1899822138089391
$ws.Range("Q9").Value = 28.495516024605
$ws.Range("R9").Value = 256.459644221445
$ws.Range("S9").Value = 0.1345923483767785
$ws.Range("T9").Value = 0.1783042686533572

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.055721
$ws.Range("H10").Value = 12.167163
$ws.Range("I10").Value = 0.9105463609091901
$ws.Range("J10").Value = 0.9385313765880943
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.048527
$ws.Range("N10").Value = 6.145581
$ws.Range("O10").Value = 0.04309744714193724
$ws.Range("P10").Value = 0.05539188977340388
$ws.Range("Q10").Value = 8.308253972967
$ws.Range("R10").Value = 74.77428575670301
$ws.Range("S10").Value = 0.03924222365956713
$ws.Range("T10").Value = 0.05198702656084872

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.055721
$ws.Range("H11").Value = 12.167163
$ws.Range("I11").Value = 0.9105463609091901
$ws.Range("J11").Value = 0.9385313765880943
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.046790666666667
$ws.Range("N11").Value = 9.140372000000001
$ws.Range("O11").Value = 0.0640991794148744
$ws.Range("P11").Value = 0.08238480272441405
$ws.Range("Q11").Value = 12.356932889404
$ws.Range("R11").Value = 111.212396004636
$ws.Range("S11").Value = 0.05836527455347915
$ws.Range("T11").Value = 0.0773207223108829

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4.055721
$ws.Range("H12").Value = 12.167163
$ws.Range("I12").Value = 0.9105463609091901
$ws.Range("J12").Value = 0.9385313765880943
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.761096
$ws.Range("N12").Value = 11.283288
$ws.Range("O12").Value = 0.0791269219569728
$ws.Range("P12").Value = 0.1016995211970309
$ws.Range("Q12").Value = 15.253956030216
$ws.Range("R12").Value = 137.285604271944
$ws.Range("S12").Value = 0.07204873083786707
$ws.Range("T12").Value = 0.09544819162739952

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.055721
$ws.Range("H13").Value = 12.167163
$ws.Range("I13").Value = 0.9105463609091901
$ws.Range("J13").Value = 0.9385313765880943
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 25.151532
$ws.Range("N13").Value = 50.303064
$ws.Range("O13").Value = 0.5291445125735434
$ws.Range("P13").Value = 0.453395989142846
$ws.Range("Q13").Value = 102.007596514572
$ws.Range("R13").Value = 612.045579087432
$ws.Range("S13").Value = 0.4818106103189072
$ws.Range("T13").Value = 0.4255263618297559
